$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 cell value from 2 to 1
$ws.Range("A2").Value = 1

# Update selection to A2
$ws.Range("A2").Select()
